$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 79
$ws1.Range("F3").Value = 659
$ws1.Range("F4").Value = 239
$ws1.Range("F6").Value = 9974
$ws1.Range("F10").Value = 5008
$ws1.Range("F11").Value = 5
$ws1.Range("F12").Value = 7
$ws1.Range("F13").Value = 178
$ws1.Range("F14").Value = 119
$ws1.Range("F18").Value = 579
$ws1.Range("F19").Value = 115
$ws1.Range("F21").Value = 9
$ws1.Range("F22").Value = 1502

# Sheet "全部类型" (sheet4): update "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 79
$ws4.Range("F4").Value = 659
$ws4.Range("F5").Value = 239
$ws4.Range("F7").Value = 9974
$ws4.Range("F11").Value = 5008
$ws4.Range("F12").Value = 5
$ws4.Range("F13").Value = 7
$ws4.Range("F14").Value = 178
$ws4.Range("F15").Value = 119
$ws4.Range("F19").Value = 579
$ws4.Range("F20").Value = 115
$ws4.Range("F22").Value = 9
$ws4.Range("F23").Value = 1502

$wb.Save()
